# Apply targeted updates to column F ("dSF") values on Sheet1, matching
# the "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 4
    4  = -2
    7  = -2
    10 = -8
    12 = -2
    13 = -4
    17 = -5
    18 = -3
    21 = -2
    27 = 7
    28 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
